# Apply the diff: insert a new "StatQuery" column after column A (becomes column B),
# shifting the former "dbExcel"/"WebExcel" headers and their values one column to the right,
# and populate the new column with a StatQuery header + matching Cypher query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; existing B (dbExcel) and C (WebExcel) shift right to C and D.
$ws.Columns.Item(2).Insert()

# New header in B1.
$ws.Range("B1").Value = "StatQuery"

# New long query text in B2, matching the wrap-text style already used by A2.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.gender IN ['MALE']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# New column B should have the same width as column A (75.81640625).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
